$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.493773341178894
$ws.Range("B1").Value = 1.741562247276306
$ws.Range("C1").Value = 1.861788153648376
$ws.Range("D1").Value = 2.137217998504639
$ws.Range("E1").Value = 2.791407108306885
